## nk-convert.xlsx: add a "Ti" sheet (after "TiN") with n/k -> e1/e2 conversion
## data, fix a rounding blip in TiN!H2, and leave the "Ti" tab active/selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- fix the tiny floating point re-round in TiN!H2 -----------------------
$ws1.Range("H2").Value = 0.70304975000000003

# --- add the new "Ti" worksheet right after "TiN" --------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Ti"

# --- header row --------------------------------------------------------
$ws2.Range("A1").Value = "wl(um)"
$ws2.Range("B1").Value = "n"
$ws2.Range("C1").Value = "k"
$ws2.Range("D1").Value = "e1"
$ws2.Range("E1").Value = "e2"
$ws2.Range("H1").Value = "e1"
$ws2.Range("K1").Value = "e2"

# --- raw n/k vs wavelength data (Ti) ------------------------------------
$data = @(
    @(0.188,1.1,1.62),
    @(0.192,1.16,1.64),
    @(0.195,1.22,1.66),
    @(0.199,1.25,1.68),
    @(0.203,1.27,1.69),
    @(0.207,1.31,1.69),
    @(0.212,1.31,1.68),
    @(0.216,1.32,1.67),
    @(0.221,1.32,1.66),
    @(0.226,1.32,1.66),
    @(0.231,1.31,1.68),
    @(0.237,1.3,1.72),
    @(0.243,1.28,1.77),
    @(0.249,1.27,1.83),
    @(0.255,1.26,1.91),
    @(0.262,1.27,1.99),
    @(0.269,1.27,2.07),
    @(0.276,1.3,2.17),
    @(0.284,1.35,2.26),
    @(0.292,1.4,2.36),
    @(0.301,1.45,2.46),
    @(0.311,1.5,2.57),
    @(0.32,1.55,2.66),
    @(0.332,1.61,2.74),
    @(0.342,1.72,2.82),
    @(0.354,1.82,2.87),
    @(0.368,1.9,2.9),
    @(0.381,1.99,2.93),
    @(0.397,2.08,2.95),
    @(0.413,2.14,2.98),
    @(0.431,2.21,3.01),
    @(0.451,2.27,3.04),
    @(0.471,2.32,3.1),
    @(0.496,2.36,3.19),
    @(0.521,2.44,3.3),
    @(0.549,2.54,3.43),
    @(0.582,2.6,3.58),
    @(0.617,2.67,3.72),
    @(0.659,2.76,3.84),
    @(0.704,2.86,3.96),
    @(0.756,3,4.01),
    @(0.821,3.21,4.01),
    @(0.892,3.29,3.96),
    @(0.984,3.35,3.97),
    @(1.088,3.5,4.02),
    @(1.216,3.62,4.15),
    @(1.393,3.67,4.37),
    @(1.61,3.69,4.7),
    @(1.937,3.51,5.19)
)

$n = $data.Count
$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$lastRow = 1 + $n          # 50 — last row that actually holds a data point
$tailRow = $lastRow + 1    # 51 — leftover styled-but-empty row (from the source sheet)

# --- formula columns: first data row is a plain formula, the rest share it -
$ws2.Range("D2").Formula = "=B2^2-C2^2"
$ws2.Range("E2").Formula = "=2*B2*C2"
$ws2.Range("G2").Formula = "=A2*10^-6"
$ws2.Range("J2").Formula = "=A2*10^-6"

$ws2.Range("D3:D" + $lastRow).Formula = "=B3^2-C3^2"
$ws2.Range("E3:E" + $lastRow).Formula = "=2*B3*C3"
$ws2.Range("G3:G" + $lastRow).Formula = "=A3*10^-6"
$ws2.Range("J3:J" + $lastRow).Formula = "=A3*10^-6"

# --- H/K "paste values" copies used by the chart --------------------------
for ($i = 2; $i -le $lastRow; $i++) {
    $ws2.Cells.Item($i, 8).Value = $ws2.Cells.Item($i, 4).Value2
    $ws2.Cells.Item($i, 11).Value = $ws2.Cells.Item($i, 5).Value2
}

# --- number formatting (scientific, matches TiN's e1/e2 columns) ----------
$ws2.Range("D2:E" + $tailRow).NumberFormat = "0.00E+00"
$ws2.Range("H2:H" + $tailRow).NumberFormat = "0.00E+00"
$ws2.Range("B" + $tailRow + ":C" + $tailRow).NumberFormat = "0.00E+00"

# --- selection / active tab: "Ti" ends up the active sheet -----------------
$ws2.Range("P7").Select()
$ws1.Range("E3").Select()
$ws2.Activate()
